# BOT; UPDATE DATA
# Appends one new day of data (2020-05-10) to the "相談件数" sheet, pushing
# the trailing footnote row down by one, and updates the sheet/workbook
# metadata (dimension, print area, pane scroll position, selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# --- Insert a new row 106 for 2020-05-10, pushing the old footnote row to 107 ---
# Row 106 currently holds the trailing "※4/8より..." footnote; inserting above it
# preserves that row's formatting/content while moving it down to row 107, and the
# freshly-inserted row 106 inherits formatting from the row above (row 105).
[void]$ws.Rows.Item(106).Insert()

$ws.Range("A106").Value = 43961
$ws.Range("B106").Value = 394
$ws.Range("C106").Value = 35779
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 7232

# --- Update the print area to cover the new last row (108, incl. footnote row 107) ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$108"
    }
}

# --- Move the active selection / scroll position down to the newly added data ---
$win = $excel.ActiveWindow
[void]$ws.Range("A106").Select()
$win.ScrollRow = 83
$win.ScrollColumn = 2
